$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 3-9 (Q1..Q7) with the new ifoCAST full-series values ---

# Row 3 (Q1)
$ws.Range("B3").Value = 0.1263387051840438
$ws.Range("C3").Value = 0.3641646616460975
$ws.Range("D3").Value = 0.3918082922282284
$ws.Range("E3").Value = 0.625945917973932
$ws.Range("F3").Value = 0.6268419021867003
$ws.Range("G3").Value = 23

# Row 4 (Q2)
$ws.Range("B4").Value = 0.7173643701084924
$ws.Range("C4").Value = 0.8116168402166231
$ws.Range("D4").Value = 4.742135511821281
$ws.Range("E4").Value = 2.177644487013727
$ws.Range("F4").Value = 2.10230446137915
$ws.Range("G4").Value = 23

# Row 5 (Q3)
$ws.Range("B5").Value = 0.3243569944388665
$ws.Range("C5").Value = 1.140585590613996
$ws.Range("D5").Value = 7.644553728315862
$ws.Range("E5").Value = 2.764878610050695
$ws.Range("F5").Value = 2.807497787378002
$ws.Range("G5").Value = 23

# Row 6 (Q4)
$ws.Range("B6").Value = 0.365601066543519
$ws.Range("C6").Value = 1.140326086931019
$ws.Range("D6").Value = 7.622357104580996
$ws.Range("E6").Value = 2.76086165980496
$ws.Range("F6").Value = 2.798050767853111
$ws.Range("G6").Value = 23

# Row 7 (Q5)
$ws.Range("B7").Value = 0.6369576950222359
$ws.Range("C7").Value = 1.773555024963742
$ws.Range("D7").Value = 12.03425596399146
$ws.Range("E7").Value = 3.469042514007498
$ws.Range("F7").Value = 3.529752246274835
$ws.Range("G7").Value = 15

# Row 8 (Q6)
$ws.Range("B8").Value = 0.5160447118266327
$ws.Range("C8").Value = 1.817060265299929
$ws.Range("D8").Value = 12.0362455981401
$ws.Range("E8").Value = 3.469329272084174
$ws.Range("F8").Value = 3.551148222788933
$ws.Range("G8").Value = 15

# Row 9 (Q7) - note F9 previously had no value at all, now it gets one
$ws.Range("B9").Value = -0.05105550657186964
$ws.Range("C9").Value = 3.541832015114859
$ws.Range("D9").Value = 28.3892791991018
$ws.Range("E9").Value = 5.32815908162489
$ws.Range("F9").Value = 5.836437872642916
$ws.Range("G9").Value = 6

# --- Add two new rows (Q8 -> row 10, Q9 -> row 11) ---
# Copy formatting from row 9's label cell (A9) so the new label cells pick up
# the same bordered/bold/centered style (cellXfs index 1), then overwrite
# the copied value with the correct label text.

$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -2.219753352431116
$ws.Range("C10").Value = 2.563348172737015
$ws.Range("D10").Value = 17.30024950301405
$ws.Range("E10").Value = 4.159356861705191
$ws.Range("F10").Value = 4.061681844980783
$ws.Range("G10").Value = 4

$ws.Range("A9").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.3277785186027308
$ws.Range("C11").Value = 0.3684111220090684
$ws.Range("D11").Value = 0.2431655120773814
$ws.Range("E11").Value = 0.4931181522489123
$ws.Range("F11").Value = 0.5210120052743136
$ws.Range("G11").Value = 2
